# "Checkout address verification added"
#
# The "Address" checkout test now passes (H2/H3: Hold -> Passed), and the
# "Qty" validation test data is refreshed: the stray A2 entry is removed,
# the numeric-qty sample in A3 becomes "Abc", and the final scenario (row 5)
# now reports a passing result instead of a failure.

$wb = $excel.ActiveWorkbook

$wsAddress = $wb.Worksheets.Item("Address")
$wsQty = $wb.Worksheets.Item("Qty")

# --- Address sheet: mark H2/H3 status as "Passed" instead of "Hold" ---
$wsAddress.Range("H2").Value = "Passed"
$wsAddress.Range("H3").Value = "Passed"

# --- Qty sheet: checkout/address verification row updates ---
# A2 ("k") cleared entirely (value + formatting), matching the target's removed cell
$wsQty.Range("A2").Clear()
# A3 value changed from "AbM" to "Abc"
$wsQty.Range("A3").Value = "Abc"
# B5 status changed from "Failed" to "Passed"
$wsQty.Range("B5").Value = "Passed"

# Update selection on Qty sheet to C10
$wsQty.Activate()
$wsQty.Range("C10").Select()

# Update the workbook's document window geometry (persisted as
# bookViews/workbookView@xWindow/yWindow/windowWidth/windowHeight)
$win = $wb.Windows.Item(1)
$win.Width = 14805
$win.Height = 8010
$win.Left = 240
$win.Top = 105
